# Generate Report for Handoff
# Adds a new handoff-report row (for file 988b7dd9-84f1-460e-b218-4e94441dc0d1.md)
# to the Overview sheet and to each language sheet (zh-cn, de-de), mirroring the
# existing row that was generated for 936d0336-fa42-4f04-a49d-58fae92759f8.md.

$wb = $excel.ActiveWorkbook

$newFileGuid = "988b7dd9-84f1-460e-b218-4e94441dc0d1"
$newHash     = "c95763e9413e001ba42da2b2e1b2c14ed8882e65"

$hyperlinkColor = 15570276   # matches the workbook's existing HyperLink style (RGB FF6495ED)
$dateFormat     = "yyyy-mm-dd HH:mm:ss"

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = "$newFileGuid.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-38-13 04:38:05"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d5023247c47fd9a56df4d3a2fc4d4c59df05ea1/e2e/$newFileGuid.md",
    "",
    "",
    "$newFileGuid.md"
)
Style-AsHyperlink $wsOverview.Range("A3")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = "$newFileGuid.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "$newFileGuid.$newHash.zh-cn.xlf"
$wsZh.Range("E3").Value = "2016-03-13 04:38:01"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"
$wsZh.Range("E3").NumberFormat = $dateFormat

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d5023247c47fd9a56df4d3a2fc4d4c59df05ea1/e2e/$newFileGuid.md",
    "",
    "",
    "$newFileGuid.md"
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d5023247c47fd9a56df4d3a2fc4d4c59df05ea1/e2e/$newFileGuid.md",
    "",
    "",
    ".md"
)
$wsZh.Hyperlinks.Add(
    $wsZh.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/efeafa9a1167cf25e9f79f9024cf56de61ce5199/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newFileGuid.$newHash.zh-cn.xlf",
    "",
    "",
    "$newFileGuid.$newHash.zh-cn.xlf"
)

Style-AsHyperlink $wsZh.Range("A3")
Style-AsHyperlink $wsZh.Range("B3")
Style-AsHyperlink $wsZh.Range("D3")
$wsZh.Range("E3").NumberFormat = $dateFormat

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = "$newFileGuid.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "$newFileGuid.$newHash.de-de.xlf"
$wsDe.Range("E3").Value = "2016-03-13 04:38:05"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"
$wsDe.Range("E3").NumberFormat = $dateFormat

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d5023247c47fd9a56df4d3a2fc4d4c59df05ea1/e2e/$newFileGuid.md",
    "",
    "",
    "$newFileGuid.md"
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d5023247c47fd9a56df4d3a2fc4d4c59df05ea1/e2e/$newFileGuid.md",
    "",
    "",
    ".md"
)
$wsDe.Hyperlinks.Add(
    $wsDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/77983f01904c8614703cff778905b7d9fe668111/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newFileGuid.$newHash.de-de.xlf",
    "",
    "",
    "$newFileGuid.$newHash.de-de.xlf"
)

Style-AsHyperlink $wsDe.Range("A3")
Style-AsHyperlink $wsDe.Range("B3")
Style-AsHyperlink $wsDe.Range("D3")
$wsDe.Range("E3").NumberFormat = $dateFormat

Write-Host "Handoff report row added for $newFileGuid.md"
